$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ----- Header row -----
$headers = @("glycan", "binding_score", "monosaccharides", "motifs", "sasa", "flexibility", "has_multi_node_motifs")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# ----- Data rows -----
$data = @(
    @("Fuc(a1-2)[Gal(a1-3)]Gal", -0.2099485278169315, "['Fuc(a1-2)', 'Gal(a1-1)', 'Gal(a1-3)']", "['Fuc(a1-2)[Gal(a1-3)]Gal']", 6.385619566891085, 1.404918654625511, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)Glc", 0.1702324631448833, "['Fuc(a1-2)', 'Gal(b1-4)', 'Gal(a1-3)']", "['Fuc(a1-2)[Gal(a1-3)]Gal']", 5.841043650482035, 1.754103373790929, $true),
    @("Fuc(a1-2)[Gal(a1-3)]Gal(b1-4)GlcNAc(b1-3)Gal(b1-4)Glc", -0.04966757525604, "['Fuc(a1-2)', 'Gal(b1-4)', 'Gal(a1-3)']", "['Fuc(a1-2)[Gal(a1-3)]Gal']", 5.948417738735455, 2.055811179438615, $true),
    @("Fuc(a1-2)[GalNAc(a1-3)]Gal(b1-3)GalNAc(b1-3)Gal(a1-4)Gal(b1-4)Glc", 4.23731479523725, "['Fuc(a1-2)', 'Gal(b1-3)', 'GalNAc(a1-3)']", "['Fuc(a1-2)[GalNAc(a1-3)]Gal']", 5.951218533364838, 6.174048051268027, $true)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

# ----- Formatting: bold font, thin box border, centered-top alignment -----
# Build the new cell style exactly once on a scratch cell far away from the
# data so the incremental style-chain that each property assignment creates
# collapses onto a single final xf record (no unused intermediate styles).
$scratch = $ws.Range("Z100")
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Borders.LineStyle = 1
$scratch.Font.Bold = $true

# Copy that single cell's resulting format onto the header row and onto the
# "glycan" column of the data rows.
$scratch.Copy()
$ws.Range("A1:G1").PasteSpecial(-4122)
$ws.Range("A2:A5").PasteSpecial(-4122)

# Remove the scratch cell's content/format again so it doesn't show up in
# the sheet's used range.
$scratch.Clear()
